# atualizei dados bibi e add
#
# The source data for faturamento_diario (daily revenue) was updated: a new
# daily record was added for day 16 of the 07/2025 period. In the sheet this
# shows up as a brand-new row inserted right after the existing "day 15"
# row (the 16th row of data, which lives at worksheet row 17), pushing every
# row from the old row 17 onward down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 17 (shifts old rows 17..107 down to 18..108).
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row with the new daily record:
#   Dia=16, total_venda=30864.92, Mês=7, Ano=2025, Período="07/2025"
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 30864.92
$ws.Cells.Item(17, 3).Value = 7
$ws.Cells.Item(17, 4).Value = 2025
$ws.Cells.Item(17, 5).Value = "07/2025"
